# Doing Updates for Financials
# Insert a new "most recent period" column before column D on sheet "L",
# shifting the existing D:K data right to E:L, then populate the new
# column D with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before D, bounded to the sheet's used rows so we
#    don't blow out the sheet dimension to the full 1,048,576 rows. -4161
#    is xlShiftToRight (a single-column Range.Insert defaults to shifting
#    cells *down*, so the direction must be given explicitly).
$ws.Range("D5:D102").Insert(-4161)

# 2. The freshly inserted column D has no formatting of its own (Excel
#    seeded it from column C by default). Pull the correct number
#    formats / fonts / alignment back from column E (which now holds
#    what used to be column D) so the new column matches its neighbours.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Populate the new column D with the newest reporting period's data.

# ---- Income Statement (rows 7-35) ----
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 14087000
$ws.Range("D9").Value = 6907000
$ws.Range("D10").Value = 7180000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 21000
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 12679000
$ws.Range("D18").Value = 1408000
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 2320000
$ws.Range("D22").Value = 574000
$ws.Range("D23").Value = 834000
$ws.Range("D24").Value = 128000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 706000
$ws.Range("D27").Value = 636000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 636000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 636000

# ---- Balance Sheet (rows 38-77) ----
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 405000
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 7960000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 48186000
$ws.Range("D48").Value = 15511000
$ws.Range("D49").Value = 665000
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = "NA"
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 78316000
$ws.Range("D57").Value = 0
$ws.Range("D58").Value = 17000
$ws.Range("D59").Value = 4225000
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 11359000
$ws.Range("D62").Value = 841000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 59798000
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 15773000
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 18518000
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement (rows 80-102) ----
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 636000
$ws.Range("D83").Value = 912000
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 4222000
$ws.Range("D91").Value = -995000
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -1444000
$ws.Range("D96").Value = -80000
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -2835000
$ws.Range("D101").Value = -10000
$ws.Range("D102").Value = -67000
